$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for the two added destination columns (new shared strings)
$ws.Range("F1").Value = "Baltimore"
$ws.Range("G1").Value = "Portland"

# Fill in the full 7x7 distance matrix. Most off-diagonal cells get the thousands-
# separator number format ("#,##0", matching the existing table style); the diagonal
# cells (distance of a city to itself = 0) and a couple of edge cells (G7, F8) keep
# the default General format, matching the source file exactly.
# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").NumberFormat = "#,##0"
$ws.Range("B2").Value = 1750
$ws.Range("C2").NumberFormat = "#,##0"
$ws.Range("C2").Value = 1750
$ws.Range("D2").NumberFormat = "#,##0"
$ws.Range("D2").Value = 1000
$ws.Range("E2").NumberFormat = "#,##0"
$ws.Range("E2").Value = 2000
$ws.Range("F2").NumberFormat = "#,##0"
$ws.Range("F2").Value = 1000
$ws.Range("G2").NumberFormat = "#,##0"
$ws.Range("G2").Value = 1000

# Row 3
$ws.Range("A3").NumberFormat = "#,##0"
$ws.Range("A3").Value = 1750
$ws.Range("B3").Value = 0
$ws.Range("C3").NumberFormat = "#,##0"
$ws.Range("C3").Value = 1500
$ws.Range("D3").NumberFormat = "#,##0"
$ws.Range("D3").Value = 1250
$ws.Range("E3").NumberFormat = "#,##0"
$ws.Range("E3").Value = 1250
$ws.Range("F3").NumberFormat = "#,##0"
$ws.Range("F3").Value = 750
$ws.Range("G3").NumberFormat = "#,##0"
$ws.Range("G3").Value = 1250

# Row 4
$ws.Range("A4").NumberFormat = "#,##0"
$ws.Range("A4").Value = 1750
$ws.Range("B4").NumberFormat = "#,##0"
$ws.Range("B4").Value = 1500
$ws.Range("C4").Value = 0
$ws.Range("D4").NumberFormat = "#,##0"
$ws.Range("D4").Value = 1000
$ws.Range("E4").NumberFormat = "#,##0"
$ws.Range("E4").Value = 1250
$ws.Range("F4").NumberFormat = "#,##0"
$ws.Range("F4").Value = 750
$ws.Range("G4").NumberFormat = "#,##0"
$ws.Range("G4").Value = 500

# Row 5
$ws.Range("A5").NumberFormat = "#,##0"
$ws.Range("A5").Value = 1000
$ws.Range("B5").NumberFormat = "#,##0"
$ws.Range("B5").Value = 1250
$ws.Range("C5").NumberFormat = "#,##0"
$ws.Range("C5").Value = 1000
$ws.Range("D5").Value = 0
$ws.Range("E5").NumberFormat = "#,##0"
$ws.Range("E5").Value = 1000
$ws.Range("F5").NumberFormat = "#,##0"
$ws.Range("F5").Value = 500
$ws.Range("G5").NumberFormat = "#,##0"
$ws.Range("G5").Value = 500

# Row 6
$ws.Range("A6").NumberFormat = "#,##0"
$ws.Range("A6").Value = 2000
$ws.Range("B6").NumberFormat = "#,##0"
$ws.Range("B6").Value = 1250
$ws.Range("C6").NumberFormat = "#,##0"
$ws.Range("C6").Value = 1250
$ws.Range("D6").NumberFormat = "#,##0"
$ws.Range("D6").Value = 1000
$ws.Range("E6").Value = 0
$ws.Range("F6").NumberFormat = "#,##0"
$ws.Range("F6").Value = 1250
$ws.Range("G6").NumberFormat = "#,##0"
$ws.Range("G6").Value = 250

# Row 7
$ws.Range("A7").NumberFormat = "#,##0"
$ws.Range("A7").Value = 1000
$ws.Range("B7").NumberFormat = "#,##0"
$ws.Range("B7").Value = 750
$ws.Range("C7").NumberFormat = "#,##0"
$ws.Range("C7").Value = 750
$ws.Range("D7").NumberFormat = "#,##0"
$ws.Range("D7").Value = 500
$ws.Range("E7").NumberFormat = "#,##0"
$ws.Range("E7").Value = 1250
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 500

# Row 8
$ws.Range("A8").NumberFormat = "#,##0"
$ws.Range("A8").Value = 1000
$ws.Range("B8").NumberFormat = "#,##0"
$ws.Range("B8").Value = 1250
$ws.Range("C8").NumberFormat = "#,##0"
$ws.Range("C8").Value = 500
$ws.Range("D8").NumberFormat = "#,##0"
$ws.Range("D8").Value = 500
$ws.Range("E8").NumberFormat = "#,##0"
$ws.Range("E8").Value = 250
$ws.Range("F8").Value = 500
$ws.Range("G8").Value = 0

# Update selection and zoom level to match the saved view state
$null = $ws.Range("E8").Select()
$excel.ActiveWindow.Zoom = 116
